$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "for the following reasons.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "for the following reasons.`rFirst of all,our society is developing rapidly,and ",
    2
)
